$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("E15").Value = "***.*"
$ws.Range("G15").Value = 2
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -22.222222222222
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 50
$ws.Range("L16").Value = -40
$ws.Range("M16").Value = -72.727272727272
$ws.Range("N16").Value = -95.454545454545
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -26.666666666666
$ws.Range("I17").Value = 9
$ws.Range("J17").Value = 17
$ws.Range("K17").Value = -47.058823529411
$ws.Range("L17").Value = -18.181818181818
$ws.Range("M17").Value = 28.571428571428
$ws.Range("N17").Value = -70.967741935483
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 6
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 3
$ws.Range("J18").Value = 4
$ws.Range("K18").Value = -25
$ws.Range("L18").Value = -40
$ws.Range("M18").Value = -85.714285714285
$ws.Range("N18").Value = -97.222222222222
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -48.717948717948
$ws.Range("I19").Value = 9
$ws.Range("J19").Value = 19
$ws.Range("K19").Value = -52.631578947368
$ws.Range("L19").Value = -10
$ws.Range("M19").Value = -50
$ws.Range("N19").Value = -40
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 11
$ws.Range("H20").Value = 57.142857142857
$ws.Range("I20").Value = 5
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -93.589743589743
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -40.74074074074
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = -31.25
$ws.Range("I21").Value = 29
$ws.Range("J21").Value = 48
$ws.Range("K21").Value = -39.583333333333
$ws.Range("L21").Value = -19.444444444444
$ws.Range("M21").Value = -53.225806451612
$ws.Range("N21").Value = -90.333333333333
$ws.Range("G22").Value = 1
$ws.Range("F23").NumberFormat = "General"
$ws.Range("F23").Value = "'0"
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -59.375
$ws.Range("F24").Value = 68
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = -25.274725274725
$ws.Range("I24").Value = 29
$ws.Range("J24").Value = 54
$ws.Range("K24").Value = -46.296296296296
$ws.Range("L24").Value = -3.333333333333
$ws.Range("M24").Value = -29.268292682926
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -35.714285714285
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 48.275862068965
$ws.Range("I25").Value = 20
$ws.Range("J25").Value = 21
$ws.Range("K25").Value = -4.761904761904
$ws.Range("L25").Value = -13.043478260869
$ws.Range("M25").Value = -9.090909090909
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Value = "'0"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("E26").Value = "***.*"
$ws.Range("G26").Value = 2
$ws.Range("C27").Value = 3
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E27").Value = 200
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 700
$ws.Range("I27").Value = 4
$ws.Range("J27").NumberFormat = "#,##0"
$ws.Range("J27").Value = 1
$ws.Range("K27").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = -20
$ws.Range("C30").NumberFormat = "General"
$ws.Range("C30").Value = "'0"
$ws.Range("J41").Value = 407
$ws.Range("K41").Value = -20.817120622568
$ws.Range("L41").Value = -32.392026578073
$ws.Range("M41").Value = -40.928882438316
$ws.Range("N41").Value = -45.074224021592
$ws.Range("J43").Value = 1097
$ws.Range("K43").Value = -58.619388909845
$ws.Range("L43").Value = -69.391741071428
$ws.Range("M43").Value = -83.555688802278
$ws.Range("N43").Value = -84.468356222568
